$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @{ A="ECs"; B="Fn1"; C="Col13a1"; D="ECs"; E=3; F=1; G=19.81532133333333; H=59.445964; I=0.02448901563336021; J=0.02448901563336021; K=3; L=1; M=0.3943663333333334; N=1.183099; O=0.393328993667368; P=0.393328993667368; Q=7.814495618048445; R=70.330460562436; S=0.009632239874974014; T=0.009632239874974012 },
  @{ A="ECs"; B="Fn1"; C="Col13a1"; D="FAPs"; E=3; F=1; G=19.81532133333333; H=59.445964; I=0.02448901563336021; J=0.02448901563336021; K=1; L=0.3333333333333333; M=0.043492; N=0.130476; O=0.04337759881273123; P=0.04337759881273123; Q=0.8618079554293334; R=7.756271598864; S=0.001062274695462602; T=0.001062274695462602 },
  @{ A="ECs"; B="Fn1"; C="Col13a1"; D="sCs"; E=3; F=1; G=19.81532133333333; H=59.445964; I=0.02448901563336021; J=0.02448901563336021; K=3; L=1; M=0.564779; N=1.694337; O=0.5632934075199009; P=0.5632934075199009; Q=11.19127736731867; R=100.721496305868; S=0.0137945010629236; T=0.0137945010629236 },
  @{ A="FAPs"; B="Fn1"; C="Col13a1"; D="ECs"; E=3; F=1; G=530.2995503333333; H=1590.898651; I=0.6553774102381563; J=0.6553774102381563; K=3; L=1; M=0.3943663333333334; N=1.183099; O=0.393328993667368; P=0.393328993667368; Q=209.1322892332721; R=1882.190603099449; S=0.2577789372412998; T=0.2577789372412998 },
  @{ A="FAPs"; B="Fn1"; C="Col13a1"; D="FAPs"; E=3; F=1; G=530.2995503333333; H=1590.898651; I=0.6553774102381563; J=0.6553774102381563; K=1; L=0.3333333333333333; M=0.043492; N=0.130476; O=0.04337759881273123; P=0.04337759881273123; Q=23.06378804309733; R=207.574092387876; S=0.02842869837223751; T=0.02842869837223751 },
  @{ A="FAPs"; B="Fn1"; C="Col13a1"; D="sCs"; E=3; F=1; G=530.2995503333333; H=1590.898651; I=0.6553774102381563; J=0.6553774102381563; K=3; L=1; M=0.564779; N=1.694337; O=0.5632934075199009; P=0.5632934075199009; Q=299.5020497377097; R=2695.518447639387; S=0.369169774624619; T=0.369169774624619 },
  @{ A="M1"; B="Fn1"; C="Col13a1"; D="ECs"; E=3; F=1; G=114.1551646666666; H=342.4654939999999; I=0.1410801048907614; J=0.1410801048907614; K=3; L=1; M=0.3943663333333334; N=1.183099; O=0.393328993667368; P=0.393328993667368; Q=45.01895372065621; R=405.1705834859059; S=0.05549089568316989; T=0.05549089568316989 },
  @{ A="M1"; B="Fn1"; C="Col13a1"; D="FAPs"; E=3; F=1; G=114.1551646666666; H=342.4654939999999; I=0.1410801048907614; J=0.1410801048907614; K=1; L=0.3333333333333333; M=0.043492; N=0.130476; O=0.04337759881273123; P=0.04337759881273123; Q=4.964836421682666; R=44.68352779514399; S=0.006119716190409488; T=0.006119716190409488 },
  @{ A="M1"; B="Fn1"; C="Col13a1"; D="sCs"; E=3; F=1; G=114.1551646666666; H=342.4654939999999; I=0.1410801048907614; J=0.1410801048907614; K=3; L=1; M=0.564779; N=1.694337; O=0.5632934075199009; P=0.5632934075199009; Q=64.47243974527532; R=580.2519577074779; S=0.07946949301718202; T=0.07946949301718202 },
  @{ A="M2"; B="Fn1"; C="Col13a1"; D="ECs"; E=3; F=1; G=42.33003733333333; H=126.990112; I=0.05231411232645103; J=0.05231411232645103; K=3; L=1; M=0.3943663333333334; N=1.183099; O=0.393328993667368; P=0.393328993667368; Q=16.69354161300978; R=150.241874517088; S=0.02057665715596463; T=0.02057665715596463 },
  @{ A="M2"; B="Fn1"; C="Col13a1"; D="FAPs"; E=3; F=1; G=42.33003733333333; H=126.990112; I=0.05231411232645103; J=0.05231411232645103; K=1; L=0.3333333333333333; M=0.043492; N=0.130476; O=0.04337759881273123; P=0.04337759881273123; Q=1.841017983701333; R=16.569161853312; S=0.00226926057674095; T=0.00226926057674095 },
  @{ A="M2"; B="Fn1"; C="Col13a1"; D="sCs"; E=3; F=1; G=42.33003733333333; H=126.990112; I=0.05231411232645103; J=0.05231411232645103; K=3; L=1; M=0.564779; N=1.694337; O=0.5632934075199009; P=0.5632934075199009; Q=23.90711615508267; R=215.164045395744; S=0.02946819459374545; T=0.02946819459374545 },
  @{ A="Neutro"; B="Fn1"; C="Col13a1"; D="ECs"; E=3; F=1; G=46.22062233333333; H=138.661867; I=0.05712234103418551; J=0.05712234103418551; K=3; L=1; M=0.3943663333333334; N=1.183099; O=0.393328993667368; P=0.393328993667368; Q=18.22785735398145; R=164.050716185833; S=0.02246787291490038; T=0.02246787291490038 },
  @{ A="Neutro"; B="Fn1"; C="Col13a1"; D="FAPs"; E=3; F=1; G=46.22062233333333; H=138.661867; I=0.05712234103418551; J=0.05712234103418551; K=1; L=0.3333333333333333; M=0.043492; N=0.130476; O=0.04337759881273123; P=0.04337759881273123; Q=2.010227306521333; R=18.092045758692; S=0.002477829992624914; T=0.002477829992624914 },
  @{ A="Neutro"; B="Fn1"; C="Col13a1"; D="sCs"; E=3; F=1; G=46.22062233333333; H=138.661867; I=0.05712234103418551; J=0.05712234103418551; K=3; L=1; M=0.564779; N=1.694337; O=0.5632934075199009; P=0.5632934075199009; Q=26.10443686079767; R=234.939931747179; S=0.03217663812666021; T=0.03217663812666021 },
  @{ A="sCs"; B="Fn1"; C="Col13a1"; D="ECs"; E=3; F=1; G=56.330706; H=168.992118; I=0.06961701587708551; J=0.06961701587708552; K=3; L=1; M=0.3943663333333334; N=1.183099; O=0.393328993667368; P=0.393328993667368; Q=22.214933979298; R=199.934405813682; S=0.02738239079705922; T=0.02738239079705922 },
  @{ A="sCs"; B="Fn1"; C="Col13a1"; D="FAPs"; E=3; F=1; G=56.330706; H=168.992118; I=0.06961701587708551; J=0.06961701587708552; K=1; L=0.3333333333333333; M=0.043492; N=0.130476; O=0.04337759881273123; P=0.04337759881273123; Q=2.449935065352; R=22.049415588168; S=0.003019818985255755; T=0.003019818985255756 },
  @{ A="sCs"; B="Fn1"; C="Col13a1"; D="sCs"; E=3; F=1; G=56.330706; H=168.992118; I=0.06961701587708551; J=0.06961701587708552; K=3; L=1; M=0.564779; N=1.694337; O=0.5632934075199009; P=0.5632934075199009; Q=31.814399803974; R=286.329598235766; S=0.03921480609477054; T=0.03921480609477054 }
)

$startRow = 2
for ($i = 0; $i -lt $rows.Count; $i++) {
  $r = $startRow + $i
  $row = $rows[$i]
  $ws.Cells.Item($r, 1).Value = $row.A
  $ws.Cells.Item($r, 2).Value = $row.B
  $ws.Cells.Item($r, 3).Value = $row.C
  $ws.Cells.Item($r, 4).Value = $row.D
  $ws.Cells.Item($r, 5).Value = $row.E
  $ws.Cells.Item($r, 6).Value = $row.F
  $ws.Cells.Item($r, 7).Value = $row.G
  $ws.Cells.Item($r, 8).Value = $row.H
  $ws.Cells.Item($r, 9).Value = $row.I
  $ws.Cells.Item($r, 10).Value = $row.J
  $ws.Cells.Item($r, 11).Value = $row.K
  $ws.Cells.Item($r, 12).Value = $row.L
  $ws.Cells.Item($r, 13).Value = $row.M
  $ws.Cells.Item($r, 14).Value = $row.N
  $ws.Cells.Item($r, 15).Value = $row.O
  $ws.Cells.Item($r, 16).Value = $row.P
  $ws.Cells.Item($r, 17).Value = $row.Q
  $ws.Cells.Item($r, 18).Value = $row.R
  $ws.Cells.Item($r, 19).Value = $row.S
  $ws.Cells.Item($r, 20).Value = $row.T
}